$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.015.92"
$ws.Range("E2").Value = "'  +6.12%  "
$ws.Range("D3").Value = "'3.681.98"
$ws.Range("E3").Value = "'  +18.72%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'620.84"
$ws.Range("E5").Value = "'  +7.74%  "
$ws.Range("D6").Value = "'182.84"
$ws.Range("E6").Value = "'  +3.21%  "
$ws.Range("D7").Value = "'3.680.00"
$ws.Range("E7").Value = "'  +18.73%  "
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E9").Value = "'  +5.87%  "
$ws.Range("E10").Value = "'  +8.29%  "
$ws.Range("E11").Value = "'  +5.31%  "
$ws.Range("D12").Value = "'0.503"
$ws.Range("E12").Value = "'  +7.64%  "
$ws.Range("D13").Value = "'40.55"
$ws.Range("E13").Value = "'  +12.38%  "
$ws.Range("E14").Value = "'  +6.65%  "
$ws.Range("D15").Value = "'4.292.68"
$ws.Range("E15").Value = "'  +18.59%  "
$ws.Range("D16").Value = "'3.679.44"
$ws.Range("E16").Value = "'  +18.60%  "
$ws.Range("D17").Value = "'71.025.96"
$ws.Range("E17").Value = "'  +6.18%  "
$ws.Range("D19").Value = "'7.55"
$ws.Range("E19").Value = "'  +7.63%  "
$ws.Range("D20").Value = "'520.52"
$ws.Range("E20").Value = "'  +8.57%  "
$ws.Range("D21").Value = "'16.94"
$ws.Range("E21").Value = "'  +1.17%  "
$ws.Range("D22").Value = "'9.28"
$ws.Range("E22").Value = "'  +19.44%  "
$ws.Range("D23").Value = "'0.744"
$ws.Range("E23").Value = "'  +7.97%  "
$ws.Range("E24").Value = "'  +13.06%  "
$ws.Range("D25").Value = "'88.63"
$ws.Range("E25").Value = "'  +6.13%  "
$ws.Range("D26").Value = "'13.54"
$ws.Range("E26").Value = "'  +8.02%  "
$ws.Range("D27").Value = "'11.05"
$ws.Range("E27").Value = "'  +9.26%  "
$ws.Range("E28").Value = "'  +0.07%  "
$ws.Range("D29").Value = "'2.55"
$ws.Range("E29").Value = "'  +10.70%  "
$ws.Range("D30").Value = "'8.21"
$ws.Range("E30").Value = "'  +3.79%  "
$ws.Range("E31").Value = "'  +12.28%  "
$ws.Range("E32").Value = "'  +18.69%  "
$ws.Range("D33").Value = "'31.72"
$ws.Range("E33").Value = "'  +13.35%  "
$ws.Range("E34").Value = "'  +4.51%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "'  -0.18%  "
$ws.Range("E36").Value = "'  +9.87%  "
$ws.Range("E37").Value = "'  +9.63%  "
$ws.Range("D38").Value = "'0.348"
$ws.Range("E38").Value = "'  +11.43%  "
$ws.Range("E39").Value = "'  +9.84%  "
$ws.Range("E40").Value = "'  +7.07%  "
$ws.Range("E41").Value = "'  +4.88%  "
$ws.Range("D42").Value = "'45.63"
$ws.Range("E42").Value = "'  -6.00%  "
$ws.Range("D43").Value = "'433.97"
$ws.Range("E43").Value = "'  +16.60%  "
$ws.Range("D44").Value = "'8.84"
$ws.Range("E44").Value = "'  +6.15%  "
$ws.Range("D45").Value = "'3.117.96"
$ws.Range("E45").Value = "'  +11.54%  "
$ws.Range("E46").Value = "'  +5.35%  "
$ws.Range("E47").Value = "'  +7.20%  "
$ws.Range("D48").Value = "'28.32"
$ws.Range("E48").Value = "'  +10.38%  "
$ws.Range("D49").Value = "'140.31"
$ws.Range("E49").Value = "'  +3.29%  "
$ws.Range("E50").Value = "'  +0.00%  "
$ws.Range("D51").Value = "'2.47"
$ws.Range("E51").Value = "'  +9.96%  "
